$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "58.254.34", European-style grouped
# numbers, or fixed-decimal strings like "5.00") that must stay TEXT, not be
# coerced to a Double by Excel auto-detection. Force text entry by temporarily
# switching the cell to a Text number format, then restore the original style
# so no visible formatting change remains.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" '58.254.34'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("E3").Value = '  +2.52%  '
$ws.Range("E4").Value = '  +0.19%  '
Set-TextValue "D5" '521.94'
$ws.Range("E5").Value = '  +0.88%  '
Set-TextValue "D6" '132.78'
$ws.Range("E6").Value = '  +0.28%  '
Set-TextValue "D7" '0.998'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.87%  '
Set-TextValue "D9" '2.517.20'
$ws.Range("E9").Value = '  +2.28%  '
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("E11").Value = '  -1.34%  '
Set-TextValue "D12" '5.17'
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("E13").Value = '  -1.04%  '
Set-TextValue "D14" '2.965.88'
$ws.Range("E14").Value = '  +2.52%  '
Set-TextValue "D15" '58.317.23'
$ws.Range("E15").Value = '  +0.86%  '
Set-TextValue "D16" '22.11'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("E17").Value = '  +0.89%  '
Set-TextValue "D18" '2.517.63'
$ws.Range("E18").Value = '  +2.30%  '
Set-TextValue "D19" '10.65'
$ws.Range("E19").Value = '  +0.93%  '
Set-TextValue "D20" '321.50'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("E21").Value = '  +0.88%  '
Set-TextValue "D22" '6.19'
$ws.Range("E22").Value = '  +8.95%  '
$ws.Range("E23").Value = '  +0.13%  '
Set-TextValue "D24" '64.61'
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("E25").Value = '  -0.25%  '
Set-TextValue "D26" '0.999'
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("E30").Value = '  +2.41%  '
Set-TextValue "D31" '168.34'
$ws.Range("E31").Value = '  +0.08%  '
Set-TextValue "D32" '1.19'
$ws.Range("E32").Value = '  +2.91%  '
Set-TextValue "D33" '6.32'
$ws.Range("E33").Value = '  +1.97%  '
Set-TextValue "D35" '0.997'
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +0.87%  '
$ws.Range("E37").Value = '  -5.86%  '
Set-TextValue "D38" '3.94'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("E39").Value = '  +1.51%  '
Set-TextValue "D40" '36.32'
$ws.Range("E40").Value = '  +0.06%  '
Set-TextValue "D41" '0.771'
$ws.Range("E41").Value = '  -1.71%  '
Set-TextValue "D42" '278.47'
$ws.Range("E42").Value = '  +2.96%  '
$ws.Range("E43").Value = '  +2.36%  '
Set-TextValue "D44" '130.52'
$ws.Range("E44").Value = '  +6.18%  '
Set-TextValue "D45" '5.00'
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("E46").Value = '  +1.82%  '
$ws.Range("E47").Value = '  +1.85%  '
$ws.Range("E48").Value = '  +3.91%  '
Set-TextValue "D49" '17.73'
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("E50").Value = '  +1.45%  '
Set-TextValue "D51" '16.87'
